# Minions_Technical_Document_RobotMakers.docx
#
# The author clicked right after "...(such as a human" and typed
# " or a minion", leaving the cursor (and therefore Word's "_GoBack"
# last-edit bookmark) positioned right after the freshly typed text -
# i.e. right before "), if there is, ...". That single edit:
#   1) splits the old run "...such as a human), if there is..."
#      into "...such as a human" + " or a minion" + "), if there is..."
#   2) moves the "_GoBack" bookmark from its old location (end of the
#      "Therefore, we created our own function." paragraph) to the new
#      caret position.

$d = $word.ActiveDocument

# --- locate the insertion point, right after "...such as a human" ---
$r = $d.Content
$found = $r.Find.Execute("(such as a human", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the target sentence fragment."
}
$r.Collapse(0)
$leftPos = $r.Start

# --- type the new text at the caret ---
$r.InsertAfter(" or a minion")
$rightPos = $r.End

# --- temporary marker so the "human" / " or a minion" runs do not get
#     coalesced back together once the document is serialised ---
$leftRange = $d.Range($leftPos, $leftPos)
$d.Bookmarks.Add("KeepSplit", $leftRange)

# --- relocate "_GoBack" (removing it from its old spot automatically,
#     since a bookmark name is unique in the document) to sit right
#     after the newly typed text, i.e. right before "), if there is" ---
$rightRange = $d.Range($rightPos, $rightPos)
$d.Bookmarks.Add("_GoBack", $rightRange)

# --- drop the temporary helper marker again; the run split it forced
#     stays intact ---
$d.Bookmarks.Item("KeepSplit").Delete()
